$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 4655.9023
$ws.Cells.Item(40, 9).Value = 1500.4286
$ws.Cells.Item(40, 10).Value = 5305.5586
$ws.Cells.Item(40, 11).Value = 1500.4286
$ws.Cells.Item(40, 12).Value = 5305.5586
$ws.Cells.Item(40, 13).Value = -1325.4286
$ws.Cells.Item(40, 14).Value = -5655.5586

$ws.Cells.Item(51, 9).Value = 2994.681
$ws.Cells.Item(51, 10).Value = 3000
$ws.Cells.Item(51, 11).Value = 2994.681
$ws.Cells.Item(51, 12).Value = 3000
$ws.Cells.Item(51, 13).Value = -2510.681
$ws.Cells.Item(51, 14).Value = -3968

$ws.Cells.Item(62, 8).Value = 4284.2
$ws.Cells.Item(62, 9).Value = 3471.1428
$ws.Cells.Item(62, 10).Value = 4995.625
$ws.Cells.Item(62, 11).Value = 3471.1428
$ws.Cells.Item(62, 12).Value = 4995.625
$ws.Cells.Item(62, 13).Value = -2847.1428
$ws.Cells.Item(62, 14).Value = -6243.625

$ws.Cells.Item(65, 8).Value = 4284.2
$ws.Cells.Item(65, 9).Value = 3471.1428
$ws.Cells.Item(65, 10).Value = 4995.625
$ws.Cells.Item(65, 11).Value = 17355.714
$ws.Cells.Item(65, 12).Value = 24978.125
$ws.Cells.Item(65, 13).Value = -14235.714
$ws.Cells.Item(65, 14).Value = -31218.125

$ws.Cells.Item(106, 8).Value = 2299.1428
$ws.Cells.Item(106, 9).Value = 2857.2
$ws.Cells.Item(106, 10).Value = 904
$ws.Cells.Item(106, 11).Value = 2857.2
$ws.Cells.Item(106, 12).Value = 904
$ws.Cells.Item(106, 13).Value = -2226.2
$ws.Cells.Item(106, 14).Value = -2166

$ws.Cells.Item(112, 8).Value = 2283.3333
$ws.Cells.Item(112, 9).Value = 1923.3334
$ws.Cells.Item(112, 10).Value = 2643.3333
$ws.Cells.Item(112, 11).Value = 5770.0002
$ws.Cells.Item(112, 12).Value = 7929.999899999999
$ws.Cells.Item(112, 13).Value = -4662.0002
$ws.Cells.Item(112, 14).Value = -10145.9999

$ws.Cells.Item(132, 8).Value = 8151.407
$ws.Cells.Item(132, 9).Value = 8447.615
$ws.Cells.Item(132, 10).Value = 450
$ws.Cells.Item(132, 11).Value = 25342.845
$ws.Cells.Item(132, 12).Value = 1350
$ws.Cells.Item(132, 13).Value = -22812.845
$ws.Cells.Item(132, 14).Value = -6410

$ws.Cells.Item(137, 8).Value = 23816192
$ws.Cells.Item(137, 9).Value = 41668904
$ws.Cells.Item(137, 10).Value = 12575.667
$ws.Cells.Item(137, 11).Value = 125006712
$ws.Cells.Item(137, 12).Value = 37727.001
$ws.Cells.Item(137, 13).Value = -125004162
$ws.Cells.Item(137, 14).Value = -42827.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1358670.8
$ws.Cells.Item(32, 9).Value = 1463178.4
$ws.Cells.Item(32, 10).Value = 34908
$ws.Cells.Item(32, 11).Value = 1463178.4
$ws.Cells.Item(32, 12).Value = 34908
$ws.Cells.Item(32, 13).Value = -1462891.4
$ws.Cells.Item(32, 14).Value = -35482

$ws.Cells.Item(61, 8).Value = 12503628
$ws.Cells.Item(61, 9).Value = 3502.5
$ws.Cells.Item(61, 10).Value = 25003754
$ws.Cells.Item(61, 11).Value = 3502.5
$ws.Cells.Item(61, 12).Value = 25003754
$ws.Cells.Item(61, 13).Value = -3290.5
$ws.Cells.Item(61, 14).Value = -25004178

$ws.Cells.Item(122, 8).Value = 2000
$ws.Cells.Item(122, 9).Value = 2000
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 6000
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -3550

$ws.Cells.Item(132, 8).Value = 3910.25
$ws.Cells.Item(132, 9).Value = 1939.6111
$ws.Cells.Item(132, 10).Value = 5522.591
$ws.Cells.Item(132, 11).Value = 5818.8333
$ws.Cells.Item(132, 12).Value = 16567.773
$ws.Cells.Item(132, 13).Value = -3288.8333
$ws.Cells.Item(132, 14).Value = -21627.773

$ws.Cells.Item(136, 8).Value = 12503628
$ws.Cells.Item(136, 9).Value = 3502.5
$ws.Cells.Item(136, 10).Value = 25003754
$ws.Cells.Item(136, 11).Value = 10507.5
$ws.Cells.Item(136, 12).Value = 75011262
$ws.Cells.Item(136, 13).Value = -7957.5
$ws.Cells.Item(136, 14).Value = -75016362

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 38226.82
$ws.Cells.Item(20, 9).Value = 44200.676
$ws.Cells.Item(20, 10).Value = 19762.182
$ws.Cells.Item(20, 11).Value = 44200.676
$ws.Cells.Item(20, 12).Value = 19762.182
$ws.Cells.Item(20, 13).Value = -43953.676
$ws.Cells.Item(20, 14).Value = -20256.182

$ws.Cells.Item(80, 8).Value = 833.93335
$ws.Cells.Item(80, 9).Value = 1118.6666
$ws.Cells.Item(80, 10).Value = 644.1111
$ws.Cells.Item(80, 11).Value = 1118.6666
$ws.Cells.Item(80, 12).Value = 644.1111
$ws.Cells.Item(80, 13).Value = -120.6666
$ws.Cells.Item(80, 14).Value = -2640.1111

$ws.Cells.Item(83, 8).Value = 833.93335
$ws.Cells.Item(83, 9).Value = 1118.6666
$ws.Cells.Item(83, 10).Value = 644.1111
$ws.Cells.Item(83, 11).Value = 5593.333000000001
$ws.Cells.Item(83, 12).Value = 3220.5555
$ws.Cells.Item(83, 13).Value = -601.3330000000005
$ws.Cells.Item(83, 14).Value = -13204.5555

$ws.Cells.Item(99, 8).Value = 14996.714
$ws.Cells.Item(99, 9).Value = 14996.714
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 11).Value = 14996.714
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 13).ClearContents()
$ws.Cells.Item(99, 14).Value = -13498.714

$ws.Cells.Item(105, 8).Value = 5913.067
$ws.Cells.Item(105, 9).Value = 2891.4167
$ws.Cells.Item(105, 10).Value = 17999.666
$ws.Cells.Item(105, 11).Value = 2891.4167
$ws.Cells.Item(105, 12).Value = 17999.666
$ws.Cells.Item(105, 13).Value = -1144.4167
$ws.Cells.Item(105, 14).Value = -21493.666

$ws.Cells.Item(109, 8).Value = 70000
$ws.Cells.Item(109, 9).Value = 0
$ws.Cells.Item(109, 10).Value = 70000
$ws.Cells.Item(109, 11).Value = 0
$ws.Cells.Item(109, 12).Value = 70000
$ws.Cells.Item(109, 14).Value = -72774

$ws.Cells.Item(132, 8).Value = 0
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 12).ClearContents()
$ws.Cells.Item(132, 14).Value = 0

$ws.Cells.Item(134, 8).Value = 3473677.2
$ws.Cells.Item(134, 9).Value = 1518.2174
$ws.Cells.Item(134, 10).Value = 83333336
$ws.Cells.Item(134, 11).Value = 4554.6522
$ws.Cells.Item(134, 12).Value = 250000008
$ws.Cells.Item(134, 13).Value = -2019.6522
$ws.Cells.Item(134, 14).Value = -250005078

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 9021974
$ws.Cells.Item(58, 9).Value = 16668137
$ws.Cells.Item(58, 10).Value = 3924532
$ws.Cells.Item(58, 11).Value = 16668137
$ws.Cells.Item(58, 12).Value = 3924532
$ws.Cells.Item(58, 13).Value = -16667934
$ws.Cells.Item(58, 14).Value = -3924938

$ws.Cells.Item(62, 8).Value = 4529.8184
$ws.Cells.Item(62, 9).Value = 4430.6665
$ws.Cells.Item(62, 10).Value = 4648.8
$ws.Cells.Item(62, 11).Value = 4430.6665
$ws.Cells.Item(62, 12).Value = 4648.8
$ws.Cells.Item(62, 13).Value = -3806.6665
$ws.Cells.Item(62, 14).Value = -5896.8

$ws.Cells.Item(65, 8).Value = 4529.8184
$ws.Cells.Item(65, 9).Value = 4430.6665
$ws.Cells.Item(65, 10).Value = 4648.8
$ws.Cells.Item(65, 11).Value = 22153.3325
$ws.Cells.Item(65, 12).Value = 23244
$ws.Cells.Item(65, 13).Value = -19033.3325
$ws.Cells.Item(65, 14).Value = -29484

$ws.Cells.Item(99, 8).Value = 1941.2142
$ws.Cells.Item(99, 9).Value = 1667.7778
$ws.Cells.Item(99, 10).Value = 2433.4
$ws.Cells.Item(99, 11).Value = 1667.7778
$ws.Cells.Item(99, 12).Value = 2433.4
$ws.Cells.Item(99, 13).Value = -169.7778000000001
$ws.Cells.Item(99, 14).Value = -5429.4

$ws.Cells.Item(119, 8).Value = 68900
$ws.Cells.Item(119, 9).Value = 0
$ws.Cells.Item(119, 10).Value = 68900
$ws.Cells.Item(119, 11).Value = 0
$ws.Cells.Item(119, 12).Value = 68900
$ws.Cells.Item(119, 14).Value = -78576

$ws.Cells.Item(126, 8).Value = 1941.2142
$ws.Cells.Item(126, 9).Value = 1667.7778
$ws.Cells.Item(126, 10).Value = 2433.4
$ws.Cells.Item(126, 11).Value = 5003.3334
$ws.Cells.Item(126, 12).Value = 7300.200000000001
$ws.Cells.Item(126, 13).Value = -2533.3334
$ws.Cells.Item(126, 14).Value = -12240.2

$ws.Cells.Item(134, 8).Value = 2925.7605
$ws.Cells.Item(134, 9).Value = 2236.3408
$ws.Cells.Item(134, 10).Value = 4049.2593
$ws.Cells.Item(134, 11).Value = 6709.0224
$ws.Cells.Item(134, 12).Value = 12147.7779
$ws.Cells.Item(134, 13).Value = -4174.0224
$ws.Cells.Item(134, 14).Value = -17217.7779

$ws.Cells.Item(136, 8).Value = 9021974
$ws.Cells.Item(136, 9).Value = 16668137
$ws.Cells.Item(136, 10).Value = 3924532
$ws.Cells.Item(136, 11).Value = 50004411
$ws.Cells.Item(136, 12).Value = 11773596
$ws.Cells.Item(136, 13).Value = -50001861
$ws.Cells.Item(136, 14).Value = -11778696

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(46, 8).Value = 18150
$ws.Cells.Item(46, 9).Value = 17500
$ws.Cells.Item(46, 10).Value = 18800
$ws.Cells.Item(46, 11).Value = 17500
$ws.Cells.Item(46, 12).Value = 18800
$ws.Cells.Item(46, 13).Value = -17344
$ws.Cells.Item(46, 14).Value = -19112

$ws.Cells.Item(80, 8).Value = 1499
$ws.Cells.Item(80, 9).Value = 1499
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 11).Value = 1499
$ws.Cells.Item(80, 12).Value = 0
$ws.Cells.Item(80, 13).Value = -501

$ws.Cells.Item(83, 8).Value = 1499
$ws.Cells.Item(83, 9).Value = 1499
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 11).Value = 7495
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 13).Value = -2503

$ws.Cells.Item(126, 8).Value = 2113.0952
$ws.Cells.Item(126, 9).Value = 1804.875
$ws.Cells.Item(126, 10).Value = 3099.4
$ws.Cells.Item(126, 11).Value = 5414.625
$ws.Cells.Item(126, 12).Value = 9298.200000000001
$ws.Cells.Item(126, 13).Value = -2944.625
$ws.Cells.Item(126, 14).Value = -14238.2

$ws.Cells.Item(132, 8).Value = 61795.8
$ws.Cells.Item(132, 9).Value = 36326.668
$ws.Cells.Item(132, 10).Value = 99999.5
$ws.Cells.Item(132, 11).Value = 108980.004
$ws.Cells.Item(132, 12).Value = 299998.5
$ws.Cells.Item(132, 13).Value = -106450.004
$ws.Cells.Item(132, 14).Value = -305058.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3370
$ws.Cells.Item(7, 9).Value = 2658
$ws.Cells.Item(7, 10).Value = 5150
$ws.Cells.Item(7, 11).Value = 2658
$ws.Cells.Item(7, 12).Value = 5150
$ws.Cells.Item(7, 13).Value = -2546
$ws.Cells.Item(7, 14).Value = -5374

$ws.Cells.Item(16, 8).Value = 874.0454999999999
$ws.Cells.Item(16, 9).Value = 807.2778
$ws.Cells.Item(16, 10).Value = 1174.5
$ws.Cells.Item(16, 11).Value = 807.2778
$ws.Cells.Item(16, 12).Value = 1174.5
$ws.Cells.Item(16, 13).Value = -637.2778
$ws.Cells.Item(16, 14).Value = -1514.5

$ws.Cells.Item(92, 8).Value = 55000
$ws.Cells.Item(92, 9).Value = 0
$ws.Cells.Item(92, 10).Value = 55000
$ws.Cells.Item(92, 11).Value = 0
$ws.Cells.Item(92, 12).Value = 55000
$ws.Cells.Item(92, 14).Value = -59992

$ws.Cells.Item(93, 8).Value = 3227.1428
$ws.Cells.Item(93, 9).Value = 2294.5
$ws.Cells.Item(93, 10).Value = 8823
$ws.Cells.Item(93, 11).Value = 2294.5
$ws.Cells.Item(93, 12).Value = 8823
$ws.Cells.Item(93, 13).Value = -1046.5
$ws.Cells.Item(93, 14).Value = -11319

$ws.Cells.Item(122, 8).Value = 3166.0476
$ws.Cells.Item(122, 9).Value = 2934.0833
$ws.Cells.Item(122, 10).Value = 3475.3333
$ws.Cells.Item(122, 11).Value = 8802.249899999999
$ws.Cells.Item(122, 12).Value = 10425.9999
$ws.Cells.Item(122, 13).Value = -6352.249899999999
$ws.Cells.Item(122, 14).Value = -15325.9999

$ws.Cells.Item(126, 8).Value = 3370
$ws.Cells.Item(126, 9).Value = 2658
$ws.Cells.Item(126, 10).Value = 5150
$ws.Cells.Item(126, 11).Value = 7974
$ws.Cells.Item(126, 12).Value = 15450
$ws.Cells.Item(126, 13).Value = -5504
$ws.Cells.Item(126, 14).Value = -20390

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 202311.4
$ws.Cells.Item(4, 9).Value = 333851.34
$ws.Cells.Item(4, 10).Value = 5001.5
$ws.Cells.Item(4, 11).Value = 333851.34
$ws.Cells.Item(4, 12).Value = 5001.5
$ws.Cells.Item(4, 13).Value = -333738.34
$ws.Cells.Item(4, 14).Value = -5227.5

$ws.Cells.Item(113, 8).Value = 1588.3939
$ws.Cells.Item(113, 9).Value = 1374.9333
$ws.Cells.Item(113, 10).Value = 1766.2778
$ws.Cells.Item(113, 11).Value = 4124.7999
$ws.Cells.Item(113, 12).Value = 5298.8334
$ws.Cells.Item(113, 13).Value = -1954.7999
$ws.Cells.Item(113, 14).Value = -9638.8334

$ws.Cells.Item(136, 8).Value = 5889406.5
$ws.Cells.Item(136, 9).Value = 3107351.2
$ws.Cells.Item(136, 10).Value = 15626600
$ws.Cells.Item(136, 11).Value = 9322053.600000001
$ws.Cells.Item(136, 12).Value = 46879800
$ws.Cells.Item(136, 13).Value = -9319503.600000001
$ws.Cells.Item(136, 14).Value = -46884900
